# Deploying to gh-pages from @ NIH-NCPI/ncpi-fhir-ig-2@b701e861ff4aea87f49ab6a6b6da8d47ed8dfde7
# Updates the generated StructureDefinition spreadsheet to reflect a refreshed
# IG build: new publication date, FHIR core version bumped back to 4.0.1 (R4),
# and the FHIR-core element definitions (ele-1 constraint text, Extension.id
# type, R4B->R4 doc link) that come along with that version change.

$wb = $excel.ActiveWorkbook

$metadata = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

# --- Metadata sheet -------------------------------------------------------
# Date
$metadata.Range("B8").Value = "2025-06-13T15:45:04+00:00"
# FHIR Version
$metadata.Range("B15").Value = "4.0.1"

# --- Elements sheet --------------------------------------------------------
# Extension row: ele-1 invariant text loses the "unless an empty Parameters
# resource ... or `$this is Parameters" clause (R4 wording, not R4B).
$elements.Range("AJ2").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}" + [char]10 + "ext-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"

# Extension.id row: Type(s) changes from "id" to "string"
$elements.Range("K3").Value = "string" + [char]10

# Extension.value[x] row: Definition text's Extensibility link moves from the
# R4B spec to the R4 spec.
$elements.Range("M6").Value = "Value of extension - must be one of a constrained set of the data types (see [Extensibility](http://hl7.org/fhir/R4/extensibility.html) for a list)."
